$d = $word.ActiveDocument

$replacements = @(
    @("86×62=5332", "71×47=3337"),
    @("31×61=1891", "97×16=1552"),
    @("79×22=1738", "79×16=1264"),
    @("59×92=5428", "42×78=3276"),
    @("49×58=2842", "94×36=3384"),
    @("27×21=567",  "29×63=1827"),
    @("40×78=3120", "88×98=8624"),
    @("45×23=1035", "19×45=855"),
    @("28×87=2436", "87×74=6438"),
    @("92×33=3036", "67×57=3819"),
    @("30×55=1650", "88×58=5104"),
    @("58×46=2668", "62×85=5270"),
    @("74×23=1702", "95×72=6840"),
    @("67×31=2077", "92×25=2300"),
    @("85×68=5780", "88×52=4576"),
    @("84×98=8232", "52×77=4004"),
    @("36×77=2772", "64×40=2560"),
    @("82×45=3690", "44×29=1276"),
    @("99×84=8316", "87×48=4176"),
    @("96×61=5856", "70×82=5740"),
    @("94×61=5734", "64×22=1408"),
    @("57×40=2280", "51×68=3468"),
    @("59×80=4720", "56×57=3192"),
    @("77×14=1078", "19×34=646"),
    @("43×92=3956", "89×69=6141")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
